$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "TOTAL"
$ws.Range("G24").Formula = "=SUM(G2:G22)"
$ws.Range("H24").Formula = "=SUM(H2:H22)"

$ws.Range("A26").Select()
